# Broodstock_survival_20181112-20190221.xlsx — "getting repo up to date"
#
# 1. Remove the stray `_xlchart.v1.*` hidden defined names left over from a
#    deleted chart (workbook.xml <definedNames>).
# 2. Append the ChiSq (prop.test) pairwise p-value summary table in rows
#    43-49 of the "Tanks1-4" sheet.
# 3. Leave the sheet scrolled/selected the way the author left it (cell
#    D46 selected) after adding the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. drop the leftover "_xlchart.v1.*" defined names -------------------
while ($wb.Names.Count -gt 0) {
    $wb.Names.Item(1).Delete()
}

# --- 2. new ChiSq test (prop.test) p-value table ---------------------------
# Column A (pairwise comparisons) written top-to-bottom first, then the two
# header cells, matching the order the values were originally entered in.
$ws.Range("A44").Value = "1_v_2"
$ws.Range("B44").Value = 0.7025
$ws.Range("A45").Value = "1_v_3"
$ws.Range("B45").Value = 0.03647
$ws.Range("A46").Value = "1_v_4"
$ws.Range("B46").Value = 0.02829
$ws.Range("A47").Value = "2_v_3"
$ws.Range("B47").Value = 0.005765
$ws.Range("A48").Value = "2_v_4"
$ws.Range("B48").Value = 0.004089
$ws.Range("A49").Value = "3_v_4"
$ws.Range("B49").Value = 1

$ws.Range("B43").Value = "p.value"
$ws.Range("A43").Value = "ChiSq test (prop.test)"

# --- 3. leave the selection where the author left it -----------------------
$ws.Range("D46").Select()
